$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Portfolio)
$ws.Range("B2").Value = 0.00056208821504521208
$ws.Range("C2").Value = 0.00471785372740519612
$ws.Range("D2").Value = 0.00002225814379319111
$ws.Range("E2").Value = 0.14164623019139341298
$ws.Range("F2").Value = 0.07489360610815957997
$ws.Range("G2").Value = 1.29845837641938000040

# Row 3 (Benchmark)
$ws.Range("B3").Value = 0.00025400425953128432
$ws.Range("C3").Value = 0.00481259734448367341
$ws.Range("D3").Value = 0.00002316109320013131
$ws.Range("E3").Value = 0.06400907340188365102
$ws.Range("F3").Value = 0.07639761440276188198
$ws.Range("G3").Value = 0.25667127900756481207
